$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "70.960.15"
$ws.Range("E2").Value = "  +2.03%  "
$ws.Range("D3").Value = "3.636.74"
$ws.Range("E3").Value = "  +3.76%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "'604.79"
$ws.Range("E5").Value = "  +0.18%  "
$ws.Range("D6").Value = "'200.80"
$ws.Range("E6").Value = "  +2.81%  "
$ws.Range("D7").Value = "'0.628"
$ws.Range("E7").Value = "  +1.22%  "
$ws.Range("D8").Value = "'1.00"
$ws.Range("E8").Value = "  +0.07%  "
$ws.Range("E9").Value = "  +9.84%  "
$ws.Range("E10").Value = "  +0.67%  "
$ws.Range("D11").Value = "'53.83"
$ws.Range("E11").Value = "  +1.45%  "
$ws.Range("E12").Value = "  +2.62%  "
$ws.Range("D13").Value = "'9.59"
$ws.Range("E13").Value = "  +1.70%  "
$ws.Range("D14").Value = "4.211.88"
$ws.Range("E14").Value = "  +3.65%  "
$ws.Range("D15").Value = "'680.79"
$ws.Range("E15").Value = "  +13.87%  "
$ws.Range("D16").Value = "71.053.71"
$ws.Range("E16").Value = "  +1.93%  "
$ws.Range("D17").Value = "'12.91"
$ws.Range("E17").Value = "  +2.26%  "
$ws.Range("D18").Value = "3.625.71"
$ws.Range("E18").Value = "  +2.52%  "
$ws.Range("D19").Value = "'19.06"
$ws.Range("E19").Value = "  +0.69%  "
$ws.Range("E20").Value = "  +0.16%  "
$ws.Range("E21").Value = "  +1.94%  "
$ws.Range("D22").Value = "'18.84"
$ws.Range("E22").Value = "  +5.45%  "
$ws.Range("E23").Value = "  +1.96%  "
$ws.Range("D24").Value = "'105.23"
$ws.Range("E24").Value = "  +3.32%  "
$ws.Range("E25").Value = "  +0.25%  "
$ws.Range("E26").Value = "  -2.05%  "
$ws.Range("D27").Value = "'10.56"
$ws.Range("E27").Value = "  -2.02%  "
$ws.Range("D28").Value = "'9.88"
$ws.Range("E28").Value = "  +4.36%  "
$ws.Range("D29").Value = "'34.36"
$ws.Range("E29").Value = "  +4.26%  "
$ws.Range("D30").Value = "'4.60"
$ws.Range("E30").Value = "  +7.82%  "
$ws.Range("E31").Value = "  +3.91%  "
$ws.Range("E32").Value = "  -0.54%  "
$ws.Range("D33").Value = "'0.115"
$ws.Range("E33").Value = "  +1.48%  "
$ws.Range("D34").Value = "'63.29"
$ws.Range("E34").Value = "  +0.32%  "
$ws.Range("D35").Value = "0.0₃0864"
$ws.Range("E35").Value = "  +7.05%  "
$ws.Range("D36").Value = "3.926.14"
$ws.Range("E36").Value = "  +4.85%  "
$ws.Range("D37").Value = "'0.999"
$ws.Range("E37").Value = "  -0.06%  "
$ws.Range("D38").Value = "'517.93"
$ws.Range("E38").Value = "  +5.13%  "
$ws.Range("D39").Value = "'3.03"
$ws.Range("E39").Value = "  -5.50%  "
$ws.Range("E40").Value = "  +0.61%  "
$ws.Range("E41").Value = "  -1.66%  "
$ws.Range("D42").Value = "'36.56"
$ws.Range("E42").Value = "  +1.51%  "
$ws.Range("E43").Value = "  +2.82%  "
$ws.Range("E44").Value = "  +2.01%  "
$ws.Range("E45").Value = "  +9.05%  "
$ws.Range("E46").Value = "  +6.86%  "
$ws.Range("E47").Value = "  +1.53%  "
$ws.Range("E48").Value = "  +3.17%  "
$ws.Range("E49").Value = "  -0.18%  "
$ws.Range("D50").Value = "'0.000247"
$ws.Range("E50").Value = "  +1.88%  "
$ws.Range("E51").Value = "  +2.73%  "
